# Fruta / hortaliza, semanal
# Apply updated weekly price/volume/date figures for
# "Hortaliza, Vega Modelo de Temuco - Ramas de apio" dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44649
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = 5000
$ws.Range("P2").Value = 5000
$ws.Range("D3").Value = 44315
$ws.Range("J3").Value = 40
$ws.Range("D4").Value = 44365
$ws.Range("J4").Value = 55
$ws.Range("D5").Value = 44504
$ws.Range("J5").Value = 55
$ws.Range("D6").Value = 44313
$ws.Range("J6").Value = 20
$ws.Range("D7").Value = 44176
$ws.Range("J7").Value = 10
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = 4000
$ws.Range("P7").Value = 4000
$ws.Range("D8").Value = 44301
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 3000
$ws.Range("P8").Value = 3000
$ws.Range("D9").Value = 44508
$ws.Range("J9").Value = 30
$ws.Range("D10").Value = 44316
$ws.Range("J10").Value = 20
$ws.Range("K10").Value = 4000
$ws.Range("L10").Value = 4000
$ws.Range("M10").Value = 4000
$ws.Range("P10").Value = 4000
$ws.Range("D11").Value = 44312
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 4000
$ws.Range("L11").Value = 4000
$ws.Range("M11").Value = 4000
$ws.Range("P11").Value = 4000
$ws.Range("D12").Value = 44680
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = 5000
$ws.Range("P12").Value = 5000
$ws.Range("D13").Value = 44497
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = 4000
$ws.Range("P13").Value = 4000
$ws.Range("D14").Value = 44280
$ws.Range("J14").Value = 55
$ws.Range("D15").Value = 44498
$ws.Range("J15").Value = 40
$ws.Range("D16").Value = 44259
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 4000
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = 4000
$ws.Range("P16").Value = 4000
$ws.Range("D17").Value = 44679
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 5000
$ws.Range("L17").Value = 5000
$ws.Range("M17").Value = 5000
$ws.Range("P17").Value = 5000
$ws.Range("D18").Value = 44509
$ws.Range("K18").Value = 4000
$ws.Range("L18").Value = 4000
$ws.Range("M18").Value = 4000
$ws.Range("P18").Value = 4000
$ws.Range("D19").Value = 44656
$ws.Range("J19").Value = 85
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = 5000
$ws.Range("P19").Value = 5000
$ws.Range("D20").Value = 44390
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = 6000
$ws.Range("P20").Value = 6000
$ws.Range("D21").Value = 44291
$ws.Range("J21").Value = 35
